$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testen F#")

# --- New citation line above the "hoofdsom/rente/looptijd" block (row 50) ---
$ws.Range("A50").Value = "Basisboek wiskunde en financiële rekenkunde, voorbeeld 4.8 p. 153"

# --- Show the formula text next to the existing A53 result (row 53) ---
$ws.Range("B53").Formula = "=FORMULATEXT(A53)"

# --- New row 55: recompute A53's result with the built-in FV function,
#     show its formula text, and note (like elsewhere in the sheet) that
#     F# reproduces it too. Grab the currency number format from an
#     existing FV cell (B23) *before* writing the formula, so the engine
#     reuses the existing style slot instead of minting a new one. ---
$ws.Range("B23").Copy()
$ws.Range("A55").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A55").Formula = "=FV(B52,C52,0,-A52,0)"
$ws.Range("B55").Formula = "=FORMULATEXT(A55)"
$ws.Range("D55").Value = "F# geeft dit ook."

# --- Keep the on-screen selection in sync with where the new content is ---
$ws.Range("D56").Select()

Write-Host "done"
